# Updates cryptos worksheet values (price/volume refresh) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '51.049.49'
$ws.Range('E2').Value = '  -1.85%  '
$ws.Range('D3').Value = '2.912.02'
$ws.Range('E3').Value = '  -2.28%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.998'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '373.85'
$ws.Range('E5').Value = '  +5.42%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '102.02'
$ws.Range('E6').Value = '  -4.98%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.542'
$ws.Range('E7').Value = '  -3.61%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.585'
$ws.Range('E9').Value = '  -4.68%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.87'
$ws.Range('E10').Value = '  -3.56%  '
$ws.Range('E11').Value = '  +0.35%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0834'
$ws.Range('E12').Value = '  -2.40%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '18.27'
$ws.Range('E13').Value = '  -5.13%  '
$ws.Range('D14').Value = '3.361.12'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.34'
$ws.Range('E15').Value = '  -3.61%  '
$ws.Range('D16').Value = '2.900.18'
$ws.Range('E16').Value = '  -2.90%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.923'
$ws.Range('E17').Value = '  -7.41%  '
$ws.Range('D18').Value = '50.917.66'
$ws.Range('E18').Value = '  -2.16%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.24'
$ws.Range('E19').Value = '  -6.36%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.20'
$ws.Range('E20').Value = '  -3.79%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.86'
$ws.Range('E21').Value = '  -5.56%  '
$ws.Range('D22').Value = '0.0₃0942'
$ws.Range('E22').Value = '  -3.22%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '68.18'
$ws.Range('E23').Value = '  -1.90%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '259.46'
$ws.Range('E24').Value = '  -1.53%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.69'
$ws.Range('E25').Value = '  -1.71%  '
$ws.Range('E26').Value = '  -5.88%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '4.09'
$ws.Range('E28').Value = '  -4.81%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '25.66'
$ws.Range('E29').Value = '  -4.45%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.08'
$ws.Range('E30').Value = '  -6.60%  '
$ws.Range('E31').Value = '  -6.40%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.27'
$ws.Range('E32').Value = '  +3.27%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '9.84'
$ws.Range('E33').Value = '  -4.42%  '
$ws.Range('E34').Value = '  -3.76%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '51.20'
$ws.Range('E35').Value = '  +0.74%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '34.10'
$ws.Range('E36').Value = '  -5.94%  '
$ws.Range('E37').Value = '  +0.22%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0422'
$ws.Range('E38').Value = '  -5.19%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.98'
$ws.Range('E39').Value = '  -6.58%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '16.97'
$ws.Range('E40').Value = '  -5.15%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.58'
$ws.Range('E41').Value = '  -4.16%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.84'
$ws.Range('E42').Value = '  -6.81%  '
$ws.Range('E43').Value = '  -3.88%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '119.24'
$ws.Range('E44').Value = '  -1.81%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '21.83'
$ws.Range('E45').Value = '  -4.05%  '
$ws.Range('E46').Value = '  -1.87%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.013.53'
$ws.Range('E47').Value = '  -5.08%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.31'
$ws.Range('E48').Value = '  -2.42%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.14'
$ws.Range('E49').Value = '  -6.42%  '
$ws.Range('D50').Value = '3.195.88'
$ws.Range('E50').Value = '  -2.31%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.236'
$ws.Range('E51').Value = '  -1.97%  '
